# Add a new "Driver License" mapping section to Sheet1, right before the
# existing "Registrant Residence Location" section (which currently starts
# at row 41). This pushes that section (and everything after it) down by
# two rows: one for the new section's single data row, one for the blank
# separator row that follows every section in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Make room: insert two rows at 41-42 (new content row + blank separator).
$ws.Rows("41:42").Insert()

# Populate the new section header/data row (string interning order matches
# the source workbook's shared-strings table: Driver License, then
# nc:IdentificationType, then Driver License ID, then the IEP path).
$ws.Range("A41").Value = "Driver License"
$ws.Range("C41").Value = "nc:IdentificationType"
$ws.Range("B41").Value = "Driver License ID"
$ws.Range("D41").Value = "nc:IdentificationID"
$ws.Range("E41").Value = "niem-xsd:string"
$ws.Range("F41").Value = "exchange:FirearmRegistrationQueryResults/nc:DriverLicense/nc:DriverLicenseIdentification/nc:IdentificationID"

# Leave the selection near the new row, matching the author's saved cursor.
$ws.Range("F46").Select()
